$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Insumos: remove the Allisson/Flores/Espinoza row (duplicate-name case
# now tracked only in NoProcesados) and surface the Eunice/Hernandez
# "sin apellido Paterno" row instead of leaving the input sheet padded
# with a trailing blank row.
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("Insumos")
$wsIns.Rows.Item(4).Delete()
$wsIns.Range("A5:C5").Copy()
$wsIns.Range("A6:C6").PasteSpecial(-4122)
$wsIns.Cells.Item(6, 1).Value = "Eunice"
$wsIns.Cells.Item(6, 3).Value = "Hernandez"
$wsIns.Cells.Item(6, 3).Value = "Hernández"
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Procesados: the two rows now hold the e-mails generated for the two
# "Ivan" entries instead of the old Aparicio / Allisson-Espinoza rows.
# ---------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("Procesados")
$wsProc.Cells.Item(2, 1).Value = "ivan.garcia@beeckerco.com"
$wsProc.Cells.Item(3, 1).Value = "ivan.hernandez@beeckerco.com"
$wsProc.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# NoProcesados: "Allisson" keeps its trailing space as typed by the
# user in the source sheet (no longer trimmed).
# ---------------------------------------------------------------------
$wsNo = $wb.Worksheets.Item("NoProcesados")
$wsNo.Cells.Item(2, 1).Value = "Allisson "

# ---------------------------------------------------------------------
# View bookkeeping: Insumos becomes the active/visible tab again.
# ---------------------------------------------------------------------
$wsIns.Select()
$wsProc.Range("A1").Select()
$wsNo.Range("A1").Select()
$wsIns.Select()
$wsIns.Range("A1").Select()
